$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 100000
$ws.Range("I8").Value = 100000
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 300000
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -299861
$ws.Range("N8").ClearContents()

$ws.Range("H28").Value = 654.4545000000001
$ws.Range("I28").Value = 275.125
$ws.Range("J28").Value = 1666
$ws.Range("K28").Value = 275.125
$ws.Range("L28").Value = 1666
$ws.Range("M28").Value = 209.875

$ws.Range("H100").Value = 1118.4445
$ws.Range("I100").Value = 633.375
$ws.Range("J100").Value = 4999
$ws.Range("K100").Value = 633.375
$ws.Range("L100").Value = 4999
$ws.Range("M100").Value = -92.375
$ws.Range("N100").Value = -6081

$ws.Range("H135").Value = 1517.421
$ws.Range("I135").Value = 645.6875
$ws.Range("J135").Value = 6166.6665
$ws.Range("K135").Value = 5811.1875
$ws.Range("L135").Value = 55499.9985
$ws.Range("M135").Value = -3276.1875
$ws.Range("N135").Value = -60569.9985

$ws.Range("H137").Value = 1763.5264
$ws.Range("I137").Value = 1532.3077
$ws.Range("J137").Value = 2264.5
$ws.Range("K137").Value = 4596.9231
$ws.Range("L137").Value = 6793.5
$ws.Range("M137").Value = -2046.9231

$ws.Range("H138").Value = 2987.6333
$ws.Range("I138").Value = 2708
$ws.Range("J138").Value = 3201.4707
$ws.Range("K138").Value = 8124
$ws.Range("L138").Value = 9604.4121
$ws.Range("M138").Value = -2984
$ws.Range("N138").Value = -19884.4121

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 125000
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 125000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 125000
$ws.Range("N24").Value = -125748

$ws.Range("H32").Value = 8534.467000000001
$ws.Range("I32").Value = 3826.5134
$ws.Range("J32").Value = 30308.75
$ws.Range("K32").Value = 3826.5134
$ws.Range("L32").Value = 30308.75
$ws.Range("M32").Value = -3539.5134

$ws.Range("H45").Value = 2337.7727
$ws.Range("I45").Value = 1562.4
$ws.Range("J45").Value = 3999.2856
$ws.Range("K45").Value = 1562.4
$ws.Range("L45").Value = 3999.2856
$ws.Range("M45").Value = -1185.4
$ws.Range("N45").Value = -4753.2856

$ws.Range("H61").Value = 4052.5
$ws.Range("I61").Value = 2176
$ws.Range("J61").Value = 5138.8945
$ws.Range("K61").Value = 2176
$ws.Range("L61").Value = 5138.8945
$ws.Range("M61").Value = -1964
$ws.Range("N61").Value = -5562.8945

$ws.Range("H74").Value = 2340.1875
$ws.Range("I74").Value = 1697.0731
$ws.Range("J74").Value = 6107
$ws.Range("K74").Value = 1697.0731
$ws.Range("L74").Value = 6107
$ws.Range("M74").Value = -823.0731000000001

$ws.Range("H77").Value = 2340.1875
$ws.Range("I77").Value = 1697.0731
$ws.Range("J77").Value = 6107
$ws.Range("K77").Value = 8485.3655
$ws.Range("L77").Value = 30535
$ws.Range("M77").Value = -4117.3655

$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()

$ws.Range("H100").Value = 125000
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 125000
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 125000
$ws.Range("N100").Value = -127164

$ws.Range("H102").Value = 1759.8
$ws.Range("I102").Value = 1759.8
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1759.8
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -137.8

$ws.Range("H122").Value = 3683.2727
$ws.Range("I122").Value = 2652.25
$ws.Range("J122").Value = 6432.6665
$ws.Range("K122").Value = 7956.75
$ws.Range("L122").Value = 19297.9995
$ws.Range("M122").Value = -5506.75

$ws.Range("H132").Value = 3006.5134
$ws.Range("I132").Value = 2937.2334
$ws.Range("J132").Value = 3303.4285
$ws.Range("K132").Value = 8811.700199999999
$ws.Range("L132").Value = 9910.2855
$ws.Range("M132").Value = -6281.700199999999
$ws.Range("N132").Value = -14970.2855

$ws.Range("H136").Value = 4052.5
$ws.Range("I136").Value = 2176
$ws.Range("J136").Value = 5138.8945
$ws.Range("K136").Value = 6528
$ws.Range("L136").Value = 15416.6835
$ws.Range("M136").Value = -3978
$ws.Range("N136").Value = -20516.6835

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4981.8
$ws.Range("I105").Value = 4727.5
$ws.Range("J105").Value = 5999
$ws.Range("K105").Value = 4727.5
$ws.Range("L105").Value = 5999
$ws.Range("M105").Value = -2980.5
$ws.Range("N105").Value = -9493

$ws.Range("H107").Value = 838.2083
$ws.Range("I107").Value = 479.33334
$ws.Range("J107").Value = 1436.3334
$ws.Range("K107").Value = 479.33334
$ws.Range("L107").Value = 1436.3334
$ws.Range("M107").Value = 1440.66666
$ws.Range("N107").Value = -5276.3334

$ws.Range("H134").Value = 6516.1763
$ws.Range("I134").Value = 3797.9092
$ws.Range("J134").Value = 11499.667
$ws.Range("K134").Value = 11393.7276
$ws.Range("L134").Value = 34499.001
$ws.Range("M134").Value = -8858.7276

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9205.967000000001
$ws.Range("I31").Value = 3334.7058
$ws.Range("J31").Value = 16883.77
$ws.Range("K31").Value = 3334.7058
$ws.Range("L31").Value = 16883.77
$ws.Range("M31").Value = -3039.7058

$ws.Range("H34").Value = 9205.967000000001
$ws.Range("I34").Value = 3334.7058
$ws.Range("J34").Value = 16883.77
$ws.Range("K34").Value = 3334.7058
$ws.Range("L34").Value = 16883.77
$ws.Range("M34").Value = -3132.7058

$ws.Range("H122").Value = 6875
$ws.Range("I122").Value = 4750
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 14250
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -11800

$ws.Range("H134").Value = 15418
$ws.Range("I134").Value = 12254.5
$ws.Range("J134").Value = 16999.75
$ws.Range("K134").Value = 36763.5
$ws.Range("L134").Value = 50999.25
$ws.Range("M134").Value = -34228.5
$ws.Range("N134").Value = -56069.25

$ws.Range("H137").Value = 38990
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 38990
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 38990
$ws.Range("N137").Value = -49190

$ws.Range("H141").Value = 247519.53
$ws.Range("I141").Value = 56799
$ws.Range("J141").Value = 316872.47
$ws.Range("K141").Value = 56799
$ws.Range("L141").Value = 316872.47
$ws.Range("M141").Value = -51619
$ws.Range("N141").Value = -327232.47

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 486.5
$ws.Range("I14").Value = 486.5
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1459.5
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -1286.5

$ws.Range("H51").Value = 861.3333
$ws.Range("I51").Value = 789.5
$ws.Range("J51").Value = 1005
$ws.Range("K51").Value = 2368.5
$ws.Range("L51").Value = 3015
$ws.Range("M51").Value = -1908.5
$ws.Range("N51").Value = -3935

$ws.Range("H86").Value = 1092.8
$ws.Range("I86").Value = 738.3333
$ws.Range("J86").Value = 1624.5
$ws.Range("K86").Value = 2214.9999
$ws.Range("L86").Value = 4873.5
$ws.Range("M86").Value = -1028.9999
$ws.Range("N86").Value = -7245.5

$ws.Range("H89").Value = 1092.8
$ws.Range("I89").Value = 738.3333
$ws.Range("J89").Value = 1624.5
$ws.Range("K89").Value = 6644.9997
$ws.Range("L89").Value = 14620.5
$ws.Range("M89").Value = -716.9997000000003
$ws.Range("N89").Value = -26476.5

$ws.Range("H104").Value = 915.5
$ws.Range("I104").Value = 916
$ws.Range("J104").Value = 915
$ws.Range("K104").Value = 2748
$ws.Range("L104").Value = 2745
$ws.Range("M104").Value = -127
$ws.Range("N104").Value = -7987

$ws.Range("H107").Value = 452.25
$ws.Range("I107").Value = 293
$ws.Range("J107").Value = 547.8
$ws.Range("K107").Value = 879
$ws.Range("L107").Value = 1643.4
$ws.Range("M107").Value = 1041
$ws.Range("N107").Value = -5483.4

$ws.Range("H137").Value = 3583.9
$ws.Range("I137").Value = 2065.8
$ws.Range("J137").Value = 5102
$ws.Range("K137").Value = 6197.400000000001
$ws.Range("L137").Value = 15306
$ws.Range("M137").Value = -1097.400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 247375.2
$ws.Range("I11").Value = 294350.3
$ws.Range("J11").Value = 12499.75
$ws.Range("K11").Value = 294350.3
$ws.Range("L11").Value = 12499.75
$ws.Range("M11").Value = -294211.3
$ws.Range("N11").Value = -12777.75

$ws.Range("H13").Value = 98000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 98000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 98000
$ws.Range("N13").Value = -98278

$ws.Range("H70").Value = 9976.556
$ws.Range("I70").Value = 9465
$ws.Range("J70").Value = 10999.667
$ws.Range("K70").Value = 9465
$ws.Range("L70").Value = 10999.667
$ws.Range("M70").Value = -9195
$ws.Range("N70").Value = -11539.667

$ws.Range("H73").Value = 9976.556
$ws.Range("I73").Value = 9465
$ws.Range("J73").Value = 10999.667
$ws.Range("K73").Value = 9465
$ws.Range("L73").Value = 10999.667
$ws.Range("M73").Value = -8529
$ws.Range("N73").Value = -12871.667

$ws.Range("H122").Value = 6643.4
$ws.Range("I122").Value = 2444
$ws.Range("J122").Value = 11442.714
$ws.Range("K122").Value = 7332
$ws.Range("L122").Value = 34328.142
$ws.Range("M122").Value = -4882

$ws.Range("H123").Value = 54997.25
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 54997.25
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 54997.25
$ws.Range("N123").Value = -59897.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 30563
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 30563
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 30563
$ws.Range("N12").Value = -30903

$ws.Range("H20").Value = 21500
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 21500
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 21500
$ws.Range("N20").Value = -21952

$ws.Range("H25").Value = 63900
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 63900
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 63900
$ws.Range("N25").Value = -64360

$ws.Range("H93").Value = 2534.3
$ws.Range("I93").Value = 1815.9166
$ws.Range("J93").Value = 3611.875
$ws.Range("K93").Value = 1815.9166
$ws.Range("L93").Value = 3611.875
$ws.Range("M93").Value = -567.9166
$ws.Range("N93").Value = -6107.875

$ws.Range("H100").Value = 3648.9697
$ws.Range("I100").Value = 3032.76
$ws.Range("J100").Value = 5574.625
$ws.Range("K100").Value = 3032.76
$ws.Range("L100").Value = 5574.625
$ws.Range("M100").Value = -2491.76
$ws.Range("N100").Value = -6656.625

$ws.Range("H122").Value = 4438.5483
$ws.Range("I122").Value = 4153.5
$ws.Range("J122").Value = 7099
$ws.Range("K122").Value = 12460.5
$ws.Range("L122").Value = 21297
$ws.Range("M122").Value = -10010.5
$ws.Range("N122").Value = -26197

$ws.Range("H136").Value = 9163.530000000001
$ws.Range("I136").Value = 5448.5835
$ws.Range("J136").Value = 9989.074000000001
$ws.Range("K136").Value = 16345.7505
$ws.Range("L136").Value = 29967.222
$ws.Range("M136").Value = -13795.7505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3372.68
$ws.Range("I126").Value = 3372.68
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 10118.04
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -7648.039999999999

$ws.Range("H136").Value = 6856.7646
$ws.Range("I136").Value = 5397.7144
$ws.Range("J136").Value = 13665.667
$ws.Range("K136").Value = 16193.1432
$ws.Range("L136").Value = 40997.001
$ws.Range("M136").Value = -13643.1432
$ws.Range("N136").Value = -46097.001
